$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.032.32'
$ws.Range('E2').Value = '  +2.47%  '
$ws.Range('D3').Value = '1.820.76'
$ws.Range('E3').Value = '  +3.08%  '
$ws.Range('D4').Value = '''1.009'
$ws.Range('E4').Value = '  +0.69%  '
$ws.Range('D5').Value = '''313.87'
$ws.Range('E5').Value = '  +2.69%  '
$ws.Range('E6').Value = '  +0.71%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('D8').Value = '''0.3700'
$ws.Range('E8').Value = '  +1.90%  '
$ws.Range('D9').Value = '''0.07279'
$ws.Range('E9').Value = '  +2.77%  '
$ws.Range('B10').Value = 'Polygon'
$ws.Range('C10').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D10').Value = '''0.8703'
$ws.Range('E10').Value = '  +2.88%  '
$ws.Range('B11').Value = 'WrappedEther'
$ws.Range('C11').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D11').Value = '2.105.65'
$ws.Range('E11').Value = '  +15.89%  '
$ws.Range('D12').Value = '''21.31'
$ws.Range('E12').Value = '  +5.18%  '
$ws.Range('D13').Value = '''5.423'
$ws.Range('D14').Value = '''6.646'
$ws.Range('E14').Value = '  +3.35%  '
$ws.Range('D15').Value = '''0.06974'
$ws.Range('E15').Value = '  +2.63%  '
$ws.Range('D16').Value = '''81.24'
$ws.Range('E16').Value = '  +2.63%  '
$ws.Range('E17').Value = '  +1.02%  '
$ws.Range('D18').Value = '''0.000008865'
$ws.Range('E18').Value = '  +2.56%  '
$ws.Range('D19').Value = '''1.008'
$ws.Range('E19').Value = '  +0.70%  '
$ws.Range('D20').Value = '''15.29'
$ws.Range('E20').Value = '  +1.94%  '
$ws.Range('D21').Value = '27.071.29'
$ws.Range('E21').Value = '  +2.59%  '
$ws.Range('D22').Value = '''5.217'
$ws.Range('E22').Value = '  +3.66%  '
$ws.Range('D23').Value = '''11.05'
$ws.Range('E23').Value = '  -1.38%  '
$ws.Range('D24').Value = '2.358.31'
$ws.Range('E24').Value = '  +18.56%  '
$ws.Range('D25').Value = '''154.56'
$ws.Range('E25').Value = '  +1.21%  '
$ws.Range('E26').Value = '  +1.79%  '
$ws.Range('D27').Value = '''18.46'
$ws.Range('E27').Value = '  +1.85%  '
$ws.Range('D28').Value = '''5.248'
$ws.Range('D29').Value = '''1.927'
$ws.Range('E29').Value = '  +13.21%  '
$ws.Range('D30').Value = '''114.98'
$ws.Range('E30').Value = '  +0.75%  '
$ws.Range('D31').Value = '''0.08993'
$ws.Range('E31').Value = '  +0.95%  '
$ws.Range('D32').Value = '''1.183'
$ws.Range('E32').Value = '  +6.52%  '
$ws.Range('E33').Value = '  +2.92%  '
$ws.Range('D34').Value = '''4.433'
$ws.Range('E34').Value = '  +2.28%  '
$ws.Range('D35').Value = '''2.818'
$ws.Range('E35').Value = '  +2.06%  '
$ws.Range('E36').Value = '  +0.66%  '
$ws.Range('D37').Value = '''1.127'
$ws.Range('E37').Value = '  +5.05%  '
$ws.Range('D38').Value = '''0.05255'
$ws.Range('E38').Value = '  +2.45%  '
$ws.Range('E39').Value = '  +2.06%  '
$ws.Range('E40').Value = '  +4.36%  '
$ws.Range('D41').Value = '''2.752'
$ws.Range('E41').Value = '  +9.43%  '
$ws.Range('D42').Value = '''0.1656'
$ws.Range('D43').Value = '''6.489'
$ws.Range('E43').Value = '  +4.07%  '
$ws.Range('D44').Value = '''8.344'
$ws.Range('E44').Value = '  +3.68%  '
$ws.Range('D45').Value = '''107.36'
$ws.Range('E45').Value = '  +2.38%  '
$ws.Range('D46').Value = '''10.41'
$ws.Range('E46').Value = '  +3.08%  '
$ws.Range('E47').Value = '  +0.84%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = '''1.656'
$ws.Range('E48').Value = '  +4.96%  '
$ws.Range('B49').Value = 'Decentraland'
$ws.Range('C49').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D49').Value = '''0.4580'
$ws.Range('E49').Value = '  +1.99%  '
$ws.Range('D50').Value = '''0.06233'
$ws.Range('E50').Value = '  +0.71%  '
$ws.Range('D51').Value = '''1.851'
$ws.Range('E51').Value = '  +7.27%  '

# Reset style for text-forced numeric-looking cells so no stray quotePrefix style sticks
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
